# Updates cryptos list price (D) and volume(1h) (E) values to match the scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text (e.g. "28.659.26") that must remain literal text,
# not be reinterpreted as a number, so force the text number format first.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.659.26"
$ws.Range("E2").Value = "  +2.38%  "
$ws.Range("D3").Value = "1.799.14"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "313.28"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").Value = "0.5366"
$ws.Range("E7").Value = "  -1.10%  "
$ws.Range("D8").Value = "0.3772"
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("D9").Value = "0.07534"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").Value = "42.53"
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("E11").Value = "  -1.66%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").Value = "20.93"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").Value = "6.170"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").Value = "7.415"
$ws.Range("E15").Value = "  +3.80%  "
$ws.Range("D16").Value = "1.794.86"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "90.37"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").Value = "0.00001063"
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("D19").Value = "0.06446"
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").Value = "17.21"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D22").Value = "5.931"
$ws.Range("E22").Value = "  -0.66%  "
$ws.Range("D23").Value = "28.649.67"
$ws.Range("E23").Value = "  +2.24%  "
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("D25").Value = "2.102"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("D26").Value = "160.39"
$ws.Range("E26").Value = "  +2.42%  "
$ws.Range("D27").Value = "20.43"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("E28").Value = "  -0.82%  "
$ws.Range("D29").Value = "1.999.35"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").Value = "122.97"
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("E31").Value = "  -3.75%  "
$ws.Range("D32").Value = "0.1029"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").Value = "5.660"
$ws.Range("E33").Value = "  -1.75%  "
$ws.Range("E34").Value = "  +2.96%  "
$ws.Range("D35").Value = "0.06500"
$ws.Range("E35").Value = "  +7.28%  "
$ws.Range("D36").Value = "0.2253"
$ws.Range("E36").Value = "  +6.80%  "
$ws.Range("D37").Value = "8.877"
$ws.Range("E37").Value = "  +2.52%  "
$ws.Range("D38").Value = "0.02307"
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("D39").Value = "5.035"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("D41").Value = "1.212"
$ws.Range("E41").Value = "  +5.01%  "
$ws.Range("D42").Value = "0.6244"
$ws.Range("E42").Value = "  -0.74%  "
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").Value = "1.395"
$ws.Range("D45").Value = "13.36"
$ws.Range("E45").Value = "  -0.42%  "
$ws.Range("D46").Value = "0.5876"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("D47").Value = "3.664"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").Value = "126.75"
$ws.Range("E48").Value = "  +3.60%  "
$ws.Range("D49").Value = "1.956"
$ws.Range("E49").Value = "  +1.58%  "
$ws.Range("D50").Value = "1.157"
$ws.Range("E50").Value = "  +1.90%  "
$ws.Range("D51").Value = "0.06897"
$ws.Range("E51").Value = "  +1.67%  "
